# Update cryptocurrency price/volume figures (Price column D, Volume(1h) column E)
# for rows 2-51 on Sheet1, matching the scraped-data refresh in the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.669.15"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -3.33%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.627.29"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.53%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.017"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +1.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.77"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.20%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5021"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -2.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.012"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.51%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2563"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06357"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.26%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.36"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07735"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.646.14"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.13%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.230"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.857.52"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.15%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5417"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.75%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅7842"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.06%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.01"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.47%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.758.28"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.03%  "

$ws.Range("E19").Value = "  +0.38%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "203.21"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -3.25%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.319"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.920"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.910"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.013"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.965"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +14.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.45"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.94%  "

$ws.Range("E27").Value = "  -1.69%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.72"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.41%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.772"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -3.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.238"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -1.99%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04988"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -4.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.240"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.23%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.176"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.528"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -3.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.339"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -1.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.643"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.8935"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.19%  "

$ws.Range("E38").Value = "  -2.05%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.123.06"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01552"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -2.13%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.579"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.24%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.012"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.622"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.35%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8086"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.88%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "99.26"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.92%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.767.29"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.28%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₈112"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.09%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4540"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +1.15%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.010"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +0.40%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "54.53"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.74%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05050"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -1.68%  "
